$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.787.93"
$ws.Range("E2").Value = "  -4.19%  "
$ws.Range("D3").Value = "2.983.75"
$ws.Range("E3").Value = "  -4.98%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "540.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.48%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.569"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.05%  "
$ws.Range("D9").Value = "2.996.39"
$ws.Range("E9").Value = "  -5.01%  "
$ws.Range("E10").Value = "  -3.77%  "
$ws.Range("E11").Value = "  -7.13%  "
$ws.Range("E12").Value = "  -3.91%  "
$ws.Range("D13").Value = "3.502.94"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.125"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.34%  "
$ws.Range("D15").Value = "61.833.28"
$ws.Range("E15").Value = "  -4.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.95"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.16%  "
$ws.Range("D17").Value = "2.985.74"
$ws.Range("E17").Value = "  -5.27%  "
$ws.Range("E18").Value = "  -5.51%  "
$ws.Range("E19").Value = "  -1.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "381.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.72%  "
$ws.Range("E22").Value = "  -5.06%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.08%  "
$ws.Range("E25").Value = "  -2.69%  "
$ws.Range("D26").Value = "3.105.73"
$ws.Range("E26").Value = "  -5.29%  "
$ws.Range("E27").Value = "  -2.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.995"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").Value = "0.0₃0942"
$ws.Range("E29").Value = "  -7.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -8.75%  "
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.50"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.52%  "
$ws.Range("E33").Value = "  -5.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "159.69"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.33%  "
$ws.Range("E35").Value = "  -5.82%  "
$ws.Range("E36").Value = "  -6.23%  "
$ws.Range("E37").Value = "  -5.35%  "
$ws.Range("E38").Value = "  -6.38%  "
$ws.Range("E39").Value = "  -8.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.58"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.62%  "
$ws.Range("D41").Value = "2.419.45"
$ws.Range("E41").Value = "  -8.25%  "
$ws.Range("E42").Value = "  -4.69%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.673"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.58%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.99"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0590"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.08%  "
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("E48").Value = "  -4.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.82"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0953"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "266.02"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.51%  "
